$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The hyperlink cell B4 displayed "Manage the money" but the link itself
# pointed at "Manage money" -- rename the visible text so it matches the
# hyperlink target (resolving the merge conflict).
$ws.Range("B4").Value = "Manage money"

# Update the remembered selection to match the new state.
$ws.Range("D9").Select()
